$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated fictitious-invoice data (rows 1-60, columns A-C).
# Each row is @(A, B, C); row 1 maps to worksheet row 1, etc.
$data = @(
    @(2104, 2515, 2611),
    @(886.08300000000008, 1925.7290000000005, 2129.4700000000003),
    @(11543.705000000004, 6488.0490000000018, 6737.002000000004),
    @(6659.7100000000019, 3682.746000000001, 4209.5439999999999),
    @(4494.0480000000007, 2247.5009999999997, 2209.3029999999999),
    @(0, 0, 0),
    @(3489, 2788, 2807),
    @(769.35000000000014, 2061.8949999999991, 2227.2549999999992),
    @(10928.412999999993, 6145.0419999999995, 6425.9440000000004),
    @(6190.8819999999969, 3433.9150000000004, 4066.0949999999993),
    @(4529.7630000000008, 2226.753999999999, 2279.1679999999992),
    @(0, 0, 0),
    @(3586, 2166, 1995),
    @(1680.8300000000006, 2365.848, 2592.3199999999993),
    @(10898.198999999999, 6496.1960000000036, 6838.2979999999998),
    @(6418.5089999999991, 3538.2750000000015, 4137.1259999999975),
    @(4330.2969999999987, 2277.2980000000002, 2393.0169999999994),
    @(0, 0, 0),
    @(1373.362499999999, 974.58799999999962, 904.50499999999965),
    @(2226.7299999999991, 1618.7199999999998, 1595.3999999999992),
    @(11257.601999999997, 5283.1869999999981, 5273.3289999999988),
    @(6446.7100000000028, 2889.52, 3212.8600000000006),
    @(4436.0550000000012, 1788.8079999999998, 1895.6379999999997),
    @(0, 0, 0),
    @(1437.6405000000002, 1305.6855000000005, 1296.9269999999997),
    @(2168.4419999999991, 1919.04, 1950.3590000000004),
    @(9275.2019999999975, 5731.1030000000019, 5661.070999999999),
    @(5620.2959999999975, 3200.3049999999994, 3448.8949999999995),
    @(3701.6649999999995, 1943.5560000000007, 1983.7219999999991),
    @(0, 0, 0),
    @(1512.3250000000003, 1711.5965000000003, 1658.5349999999992),
    @(1903.3499999999999, 1415.1059999999993, 1425.0379999999998),
    @(16518.094000000001, 10803.777, 10427.597999999998),
    @(5041.8729999999978, 3122.7190000000005, 3308.4050000000007),
    @(4795.8240000000023, 2852.4750000000004, 2852.4190000000003),
    @(5.8399999999996908, 133.69000000000005, 121.76000000000067),
    @(5059.6164999999946, 3905.7699999999977, 3782.1229999999991),
    @(1969.172, 892.43800000000022, 945.06500000000028),
    @(26598.50299999999, 14668.986999999999, 14325.897999999994),
    @(5868.4619999999932, 3155.5130000000022, 3470.2359999999999),
    @(8294.1510000000035, 4466.217999999998, 4292.9000000000015),
    @(0, 0, 0),
    @(5346.5099999999984, 4057.3584999999989, 4199.5489999999991),
    @(1541.1400000000006, 774.31000000000063, 781.94999999999982),
    @(23306.753000000026, 14927.280000000004, 14929.637999999992),
    @(5499.7400000000034, 3429.6760000000004, 3680.5840000000017),
    @(7589.8439999999946, 4775.7479999999987, 4709.6369999999997),
    @(0, 0, 0),
    @(3887.4975000000031, 2863.0859999999971, 3005.5939999999991),
    @(2235.9100000000017, 1762.1320000000005, 1902.9730000000004),
    @(16759.309999999994, 9248.1730000000025, 9299.146999999999),
    @(4869.9800000000023, 2704.0450000000005, 2941.7550000000001),
    @(5089.1450000000059, 2495.0159999999996, 2607.3530000000005),
    @(0, 0, 0),
    @(3381.9580000000005, 2380.5339999999983, 2428.5495000000019),
    @(2284.5360000000001, 1891.0830000000001, 2068.4569999999999),
    @(12252.739999999998, 6782.1049999999932, 6758.3839999999964),
    @(4328.6750000000002, 2233.5149999999994, 2605.3999999999987),
    @(3628.4829999999993, 1841.2869999999996, 1956.3139999999999),
    @(0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value2 = $data[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $data[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $data[$i][2]
}

# Narrow columns A and B to match column C's width (9.7109375).
# ColumnWidth is expressed in character units and is internally quantized
# to 1/6-character steps, so 8.8 is the closest input that reproduces the
# same stored width (9.666...) that column C already has on disk.
$ws.Columns.Item(1).ColumnWidth = 8.8
$ws.Columns.Item(2).ColumnWidth = 8.8
